$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.021.53"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "1.630.17"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").Value = "'214.11"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("D8").Value = "'0.248"
$ws.Range("E8").Value = "  -2.92%  "
$ws.Range("E9").Value = "  -3.35%  "
$ws.Range("D10").Value = "'18.24"
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").Value = "1.856.16"
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("D13").Value = "1.647.10"
$ws.Range("E13").Value = "  -4.50%  "
$ws.Range("E14").Value = "  -2.46%  "
$ws.Range("D15").Value = "'0.522"
$ws.Range("E15").Value = "  -4.06%  "
$ws.Range("D16").Value = "25.993.26"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").Value = "0.0₃0739"
$ws.Range("E17").Value = "  -3.30%  "
$ws.Range("D18").Value = "'61.26"
$ws.Range("E18").Value = "  -3.34%  "
$ws.Range("D20").Value = "'189.63"
$ws.Range("E20").Value = "  -2.93%  "
$ws.Range("E21").Value = "  -3.06%  "
$ws.Range("D22").Value = "'9.54"
$ws.Range("E22").Value = "  -3.93%  "
$ws.Range("D23").Value = "'6.05"
$ws.Range("E23").Value = "  -2.85%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").Value = "'143.77"
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("E26").Value = "  -1.23%  "
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("D28").Value = "'6.73"
$ws.Range("E28").Value = "  -2.47%  "
$ws.Range("D29").Value = "'15.11"
$ws.Range("E29").Value = "  -2.85%  "
$ws.Range("E30").Value = "  -1.54%  "
$ws.Range("E31").Value = "  -3.47%  "
$ws.Range("D32").Value = "'3.12"
$ws.Range("E32").Value = "  -4.34%  "
$ws.Range("D33").Value = "'3.11"
$ws.Range("E33").Value = "  -5.60%  "
$ws.Range("D34").Value = "'2.40"
$ws.Range("E34").Value = "  -2.27%  "
$ws.Range("D35").Value = "'1.48"
$ws.Range("E35").Value = "  -3.43%  "
$ws.Range("D36").Value = "1.131.80"
$ws.Range("E36").Value = "  -0.26%  "
$ws.Range("D37").Value = "'0.851"
$ws.Range("E37").Value = "  -6.12%  "
$ws.Range("E38").Value = "  -1.13%  "
$ws.Range("D39").Value = "'0.515"
$ws.Range("E39").Value = "  -4.81%  "
$ws.Range("E40").Value = "  -1.88%  "
$ws.Range("D41").Value = "'98.06"
$ws.Range("E41").Value = "  -1.37%  "
$ws.Range("D42").Value = "'0.772"
$ws.Range("E42").Value = "  -3.23%  "
$ws.Range("D43").Value = "1.767.03"
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("D44").Value = "'5.19"
$ws.Range("E44").Value = "  -5.52%  "
$ws.Range("D46").Value = "'54.62"
$ws.Range("E46").Value = "  -3.80%  "
$ws.Range("D47").Value = "'0.0527"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").Value = "'1.47"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.51"
$ws.Range("E50").Value = "  -3.55%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'1.01"
$ws.Range("E51").Value = "  +0.48%  "
